# Adds a new introductory paragraph ("This is really cool stuff. I think "
# followed by a Wingdings smiley symbol) at the very top of the document,
# and relocates the "_GoBack" bookmark from the end of the last paragraph
# onto the newly inserted paragraph.

$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the last paragraph.
# Remove it from there first -- it will be re-created on the new first
# paragraph below (mirrors how Word moves "_GoBack" to the most recent edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Build the new paragraph as a FlatOPC Word package fragment so that we can
# emit a genuine <w:sym> run (a plain InsertSymbol/Text assignment can't
# produce a <w:sym> element), plus the relocated bookmark, all in one shot.
$flatOpc = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This is really cool stuff. I think </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F04A"/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$startOfDoc = $d.Range(0, 0)
[void]$startOfDoc.InsertXML($flatOpc)
